$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-23 Thursday" "2025-01-24 Friday"

Replace-Text "228÷9=" "419÷2="
Replace-Text "158÷7=" "662÷9="
Replace-Text "831÷9=" "188÷2="
Replace-Text "334÷7=" "984÷7="
Replace-Text "739÷4=" "710÷5="
Replace-Text "970÷6=" "923÷3="
Replace-Text "182÷6=" "706÷6="
Replace-Text "535÷5=" "892÷9="
Replace-Text "315÷4=" "482÷5="
Replace-Text "870÷8=" "549÷3="
Replace-Text "150÷3=" "306÷7="
Replace-Text "223÷5=" "448÷7="
Replace-Text "619÷5=" "401÷9="
Replace-Text "755÷4=" "421÷2="
Replace-Text "423÷9=" "427÷3="
Replace-Text "902÷9=" "209÷9="
Replace-Text "496÷5=" "207÷4="
Replace-Text "671÷9=" "181÷6="
Replace-Text "290÷8=" "673÷2="
Replace-Text "673÷8=" "425÷8="
Replace-Text "271÷9=" "780÷5="
Replace-Text "564÷2=" "325÷6="
Replace-Text "918÷5=" "443÷9="
Replace-Text "324÷9=" "653÷7="
Replace-Text "895÷3=" "811÷2="
